$wb = $excel.ActiveWorkbook

# --- Update timestamps on the "data" sheet (column F) ---
$dataSheet = $wb.Worksheets.Item("data")
$dataSheet.Range("F2").Value = "2021-10-05 14:34:04.116187"
$dataSheet.Range("F3").Value = "2021-10-05 14:34:04.116194"
$dataSheet.Range("F4").Value = "2021-10-05 14:34:04.116198"
$dataSheet.Range("F5").Value = "2021-10-05 14:34:04.116200"
$dataSheet.Range("F6").Value = "2021-10-05 14:34:04.116203"
$dataSheet.Range("F7").Value = "2021-10-05 14:34:04.116206"
$dataSheet.Range("F8").Value = "2021-10-05 14:34:04.116208"
$dataSheet.Range("F9").Value = "2021-10-05 14:34:04.116211"
$dataSheet.Range("F10").Value = "2021-10-05 14:34:04.116213"
$dataSheet.Range("F11").Value = "2021-10-05 14:34:04.116216"
$dataSheet.Range("F12").Value = "2021-10-05 14:34:04.116219"
$dataSheet.Range("F13").Value = "2021-10-05 14:34:04.116221"

# --- Add a new "metadata" worksheet ---
$metaSheet = $wb.Worksheets.Add()
$metaSheet.Name = "metadata"

# Re-resolve the "data" sheet reference now that a new sheet has been added,
# since previously captured references can become stale/shifted.
$dataSheet = $wb.Worksheets.Item("data")

# Match page setup / outline settings used on the "data" sheet
$metaSheet.Outline.SummaryRow = 1
$metaSheet.Outline.SummaryColumn = 1

$metaSheet.PageSetup.LeftMargin = 54
$metaSheet.PageSetup.RightMargin = 54
$metaSheet.PageSetup.TopMargin = 72
$metaSheet.PageSetup.BottomMargin = 72
$metaSheet.PageSetup.HeaderMargin = 36
$metaSheet.PageSetup.FooterMargin = 36

# Header row (row 1)
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row (row 2)
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Hypercalcaemia"
$metaSheet.Range("C2").Value = 117
$metaSheet.Range("D2").Value = "'1.0"
$metaSheet.Range("D2").Style = "Normal"
$metaSheet.Range("E2").Value = "2021-04-06T02:30:49.731747Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:34:04.112992"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/117/?format=json"

# Apply the same header style used on the "data" sheet (row1 header cells /
# index column) by copying formats from an existing styled cell, so that the
# existing style (s="1") gets reused instead of a brand-new style created.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)

$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

# Move "metadata" to be positioned right after "data" (re-resolve the data
# sheet reference fresh, right before moving).
$dataSheet = $wb.Worksheets.Item("data")
$metaSheet.Move($null, $dataSheet)
